$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5906356666666667
$ws.Range("H2").Value = 1.771907
$ws.Range("I2").Value = 0.294823169192623
$ws.Range("J2").Value = 0.294823169192623
$ws.Range("M2").Value = 2.769264333333334
$ws.Range("N2").Value = 8.307793
$ws.Range("O2").Value = 0.1388016358751757
$ws.Range("P2").Value = 0.1388016358751757
$ws.Range("Q2").Value = 1.635626285694556
$ws.Range("R2").Value = 14.720636571251
$ws.Range("S2").Value = 0.04092193817783978
$ws.Range("T2").Value = 0.04092193817783978

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5906356666666667
$ws.Range("H3").Value = 1.771907
$ws.Range("I3").Value = 0.294823169192623
$ws.Range("J3").Value = 0.294823169192623
$ws.Range("M3").Value = 1.484487666666667
$ws.Range("N3").Value = 4.453463
$ws.Range("O3").Value = 0.07440579582442265
$ws.Range("P3").Value = 0.07440579582442265
$ws.Range("Q3").Value = 0.8767913626601112
$ws.Range("R3").Value = 7.891122263941
$ws.Range("S3").Value = 0.02193655253125552
$ws.Range("T3").Value = 0.02193655253125552

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5906356666666667
$ws.Range("H4").Value = 1.771907
$ws.Range("I4").Value = 0.294823169192623
$ws.Range("J4").Value = 0.294823169192623
$ws.Range("M4").Value = 0.1509683333333333
$ws.Range("N4").Value = 0.452905
$ws.Range("O4").Value = 0.007566865820567083
$ws.Range("P4").Value = 0.007566865820567083
$ws.Range("Q4").Value = 0.08916728220388889
$ws.Range("R4").Value = 0.8025055398349999
$ws.Range("S4").Value = 0.002230887362074925
$ws.Range("T4").Value = 0.002230887362074925

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.5906356666666667
$ws.Range("H5").Value = 1.771907
$ws.Range("I5").Value = 0.294823169192623
$ws.Range("J5").Value = 0.294823169192623
$ws.Range("M5").Value = 15.54651666666667
$ws.Range("N5").Value = 46.63955
$ws.Range("O5").Value = 0.7792257024798346
$ws.Range("P5").Value = 0.7792257024798346
$ws.Range("Q5").Value = 9.182327235761111
$ws.Range("R5").Value = 82.64094512185
$ws.Range("S5").Value = 0.2297337911214528
$ws.Range("T5").Value = 0.2297337911214528

$ws.Range("I6").Value = 0.2901829546991739
$ws.Range("J6").Value = 0.2901829546991739
$ws.Range("M6").Value = 2.769264333333334
$ws.Range("N6").Value = 8.307793
$ws.Range("O6").Value = 0.1388016358751757
$ws.Range("P6").Value = 0.1388016358751757
$ws.Range("Q6").Value = 1.609883204451889
$ws.Range("R6").Value = 14.488948840067
$ws.Range("S6").Value = 0.04027786881533735
$ws.Range("T6").Value = 0.04027786881533736

$ws.Range("I7").Value = 0.2901829546991739
$ws.Range("J7").Value = 0.2901829546991739
$ws.Range("M7").Value = 1.484487666666667
$ws.Range("N7").Value = 4.453463
$ws.Range("O7").Value = 0.07440579582442265
$ws.Range("P7").Value = 0.07440579582442265
$ws.Range("Q7").Value = 0.8629915653107778
$ws.Range("R7").Value = 7.766924087797001
$ws.Range("S7").Value = 0.02159129367907442
$ws.Range("T7").Value = 0.02159129367907442

$ws.Range("I8").Value = 0.2901829546991739
$ws.Range("J8").Value = 0.2901829546991739
$ws.Range("M8").Value = 0.1509683333333333
$ws.Range("N8").Value = 0.452905
$ws.Range("O8").Value = 0.007566865820567083
$ws.Range("P8").Value = 0.007566865820567083
$ws.Range("Q8").Value = 0.08776388057722223
$ws.Range("R8").Value = 0.7898749251950001
$ws.Range("S8").Value = 0.002195775481624345
$ws.Range("T8").Value = 0.002195775481624346

$ws.Range("I9").Value = 0.2901829546991739
$ws.Range("J9").Value = 0.2901829546991739
$ws.Range("M9").Value = 15.54651666666667
$ws.Range("N9").Value = 46.63955
$ws.Range("O9").Value = 0.7792257024798346
$ws.Range("P9").Value = 0.7792257024798346
$ws.Range("Q9").Value = 9.037806816827779
$ws.Range("R9").Value = 81.34026135145001
$ws.Range("S9").Value = 0.2261180167231378
$ws.Range("T9").Value = 0.2261180167231378

$ws.Range("G10").Value = 0.817256
$ws.Range("H10").Value = 2.451768
$ws.Range("I10").Value = 0.4079435387325965
$ws.Range("J10").Value = 0.4079435387325965
$ws.Range("M10").Value = 2.769264333333334
$ws.Range("N10").Value = 8.307793
$ws.Range("O10").Value = 0.1388016358751757
$ws.Range("P10").Value = 0.1388016358751757
$ws.Range("Q10").Value = 2.263197892002667
$ws.Range("R10").Value = 20.368781028024
$ws.Range("S10").Value = 0.05662323052079251
$ws.Range("T10").Value = 0.05662323052079251

$ws.Range("G11").Value = 0.817256
$ws.Range("H11").Value = 2.451768
$ws.Range("I11").Value = 0.4079435387325965
$ws.Range("J11").Value = 0.4079435387325965
$ws.Range("M11").Value = 1.484487666666667
$ws.Range("N11").Value = 4.453463
$ws.Range("O11").Value = 0.07440579582442265
$ws.Range("P11").Value = 0.07440579582442265
$ws.Range("Q11").Value = 1.213206452509333
$ws.Range("R11").Value = 10.918858072584
$ws.Range("S11").Value = 0.03035336365083003
$ws.Range("T11").Value = 0.03035336365083003

$ws.Range("G12").Value = 0.817256
$ws.Range("H12").Value = 2.451768
$ws.Range("I12").Value = 0.4079435387325965
$ws.Range("J12").Value = 0.4079435387325965
$ws.Range("M12").Value = 0.1509683333333333
$ws.Range("N12").Value = 0.452905
$ws.Range("O12").Value = 0.007566865820567083
$ws.Range("P12").Value = 0.007566865820567083
$ws.Range("Q12").Value = 0.1233797762266667
$ws.Range("R12").Value = 1.11041798604
$ws.Range("S12").Value = 0.003086854019956868
$ws.Range("T12").Value = 0.003086854019956868

$ws.Range("G13").Value = 0.817256
$ws.Range("H13").Value = 2.451768
$ws.Range("I13").Value = 0.4079435387325965
$ws.Range("J13").Value = 0.4079435387325965
$ws.Range("M13").Value = 15.54651666666667
$ws.Range("N13").Value = 46.63955
$ws.Range("O13").Value = 0.7792257024798346
$ws.Range("P13").Value = 0.7792257024798346
$ws.Range("Q13").Value = 12.70548402493333
$ws.Range("R13").Value = 114.3493562244
$ws.Range("S13").Value = 0.3178800905410171
$ws.Range("T13").Value = 0.3178800905410171

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.01412433333333333
$ws.Range("H14").Value = 0.042373
$ws.Range("I14").Value = 0.007050337375606629
$ws.Range("J14").Value = 0.007050337375606628
$ws.Range("M14").Value = 2.769264333333334
$ws.Range("N14").Value = 8.307793
$ws.Range("O14").Value = 0.1388016358751757
$ws.Range("P14").Value = 0.1388016358751757
$ws.Range("Q14").Value = 0.03911401253211112
$ws.Range("R14").Value = 0.352026112789
$ws.Range("S14").Value = 0.0009785983612060935
$ws.Range("T14").Value = 0.0009785983612060933

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.01412433333333333
$ws.Range("H15").Value = 0.042373
$ws.Range("I15").Value = 0.007050337375606629
$ws.Range("J15").Value = 0.007050337375606628
$ws.Range("M15").Value = 1.484487666666667
$ws.Range("N15").Value = 4.453463
$ws.Range("O15").Value = 0.07440579582442265
$ws.Range("P15").Value = 0.07440579582442265
$ws.Range("Q15").Value = 0.02096739863322222
$ws.Range("R15").Value = 0.188706587699
$ws.Range("S15").Value = 0.0005245859632626827
$ws.Range("T15").Value = 0.0005245859632626826

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.01412433333333333
$ws.Range("H16").Value = 0.042373
$ws.Range("I16").Value = 0.007050337375606629
$ws.Range("J16").Value = 0.007050337375606628
$ws.Range("M16").Value = 0.1509683333333333
$ws.Range("N16").Value = 0.452905
$ws.Range("O16").Value = 0.007566865820567083
$ws.Range("P16").Value = 0.007566865820567083
$ws.Range("Q16").Value = 0.002132327062777778
$ws.Range("R16").Value = 0.019190943565
$ws.Range("S16").Value = 0.00005334895691094442
$ws.Range("T16").Value = 0.00005334895691094442

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.01412433333333333
$ws.Range("H17").Value = 0.042373
$ws.Range("I17").Value = 0.007050337375606629
$ws.Range("J17").Value = 0.007050337375606628
$ws.Range("M17").Value = 15.54651666666667
$ws.Range("N17").Value = 46.63955
$ws.Range("O17").Value = 0.7792257024798346
$ws.Range("P17").Value = 0.7792257024798346
$ws.Range("Q17").Value = 0.2195841835694556
$ws.Range("R17").Value = 1.97625765215
$ws.Range("S17").Value = 0.005493804094226908
$ws.Range("T17").Value = 0.005493804094226908
